# Applies the changes described by the diff:
#  - "#unique=true" tag text changed to "#match=unique" throughout
#  - New rows/tests added to the "#convert" sheet (rows 7,8,10,11)
#  - New rows/tests added to the "#export" sheet (rows 4,5,6)
#  - Selections / active cells updated to reflect the new last-used cell

$wb = $excel.ActiveWorkbook

$wsConvert = $wb.Worksheets.Item("#convert")
$wsTagging = $wb.Worksheets.Item("#tagging")
$wsExport  = $wb.Worksheets.Item("#export")

# ---------------------------------------------------------------------
# "#convert" sheet ("sheet1"): rename the #unique=true tag to
# #match=unique on the two existing test blocks, and append two new
# test blocks (rows 7-8 and rows 10-11) that exercise the same
# match=unique / exact-comparison behavior with new data.
# ---------------------------------------------------------------------
# Introduce the renamed tag text first (so it reuses/leads the shared
# string table in the same relative order the source workbook uses).
$wsConvert.Range("D1").Value = "#match=unique"
$wsConvert.Range("D4").Value = "#match=unique"

# Introduce the brand-new unique strings in the same order they first
# appear in the target workbook: qwer, #measurement.formula.value,
# zxcv, ghjk.
$wsConvert.Range("B8").Value = "qwer"
$wsConvert.Range("C8").Value = "qwer"

$wsConvert.Range("A7").Value = "#tags"
$wsConvert.Range("B7").Value = "#measurement.formula.value"
$wsConvert.Range("C7").Value = "#measurement.assignment.assign"
$wsConvert.Range("D7").Value = "#match=unique"
$wsConvert.Range("E7").Value = "#comparison=exact"

$wsConvert.Range("C11").Value = "zxcv"
$wsConvert.Range("B11").Value = "ghjk"

$wsConvert.Range("A10").Value = "#tags"
$wsConvert.Range("B10").Value = "#measurement.formula.value"
$wsConvert.Range("C10").Value = "#measurement.assignment.assign"
$wsConvert.Range("D10").Value = "#match=unique"
$wsConvert.Range("E10").Value = "#comparison=exact"

$wsConvert.Columns.Item(4).ColumnWidth = 12.666666666666666

$wsConvert.Activate()
$wsConvert.Range("C12").Select() | Out-Null

# ---------------------------------------------------------------------
# "#tagging" sheet (sheet2): no content actually changes here (the
# shared-string table simply gets renumbered upstream), but rewriting
# the cell keeps it consistent/explicit.
# ---------------------------------------------------------------------
$wsTagging.Range("D11").Value = "*#measurement.compound"

# ---------------------------------------------------------------------
# "#export" sheet (sheet3): duplicate the two existing data rows as
# three new rows (4, 5, 6) that cover additional tracking/untracking
# test cases, re-using the same numeric measurements but with new
# identifying text in columns A and B.
# ---------------------------------------------------------------------
$wsExport.Range("A4").Value = "asdf"
$wsExport.Range("B4").Value = "qwer"
$wsExport.Range("C4").Value = 1
$wsExport.Range("D4").Value = "01_A0_Colon_T03-2017_naive_170427_UKy_GCB_rep1-quench"
$wsExport.Range("E4").Value = 289287.73343735602
$wsExport.Range("F4").Value = 0
$wsExport.Range("G4").Value = 8490014.3650100008
$wsExport.Range("H4").Value = 0
$wsExport.Range("I4").Value = 439597.55237699999
$wsExport.Range("J4").Value = "NA"
$wsExport.Range("K4").Value = 0
$wsExport.Range("L4").Value = 0
$wsExport.Range("M4").Value = 20
$wsExport.Range("N4").Value = 10
$wsExport.Range("O4").Value = 0.618176844244679
$wsExport.Range("P4").Value = 0.255757329816374
$wsExport.Range("Q4").Value = 0
$wsExport.Range("R4").Value = 0
$wsExport.Range("S4").Value = "Compound: name of assigned metabolite, noStd means assigment was NOT verified with standard compound"

$wsExport.Range("A5").Value = "zxcv"
$wsExport.Range("B5").Value = "ghjk"
$wsExport.Range("C5").Value = 0
$wsExport.Range("D5").Value = "01_A0_Colon_T03-2017_naive_170427_UKy_GCB_rep1-quench"
$wsExport.Range("E5").Value = 7989221.8338638796
$wsExport.Range("F5").Value = 8447352.8921099994
$wsExport.Range("G5").Value = 8490014.3650100008
$wsExport.Range("H5").Value = 8447352.8921099994
$wsExport.Range("I5").Value = 7839899.2880199999
$wsExport.Range("J5").Value = "NA"
$wsExport.Range("K5").Value = 0
$wsExport.Range("L5").Value = 0
$wsExport.Range("M5").Value = 20
$wsExport.Range("N5").Value = 10
$wsExport.Range("O5").Value = 0.618176844244679
$wsExport.Range("P5").Value = 0.255757329816374
$wsExport.Range("Q5").Value = 0
$wsExport.Range("R5").Value = 0
$wsExport.Range("S5").Value = "Legend"

$wsExport.Range("A6").Value = "zxcv"
$wsExport.Range("B6").Value = "ghjk"
$wsExport.Range("C6").Value = 1
$wsExport.Range("D6").Value = "01_A0_Colon_T03-2017_naive_170427_UKy_GCB_rep1-quench"
$wsExport.Range("E6").Value = 289287.73343735602
$wsExport.Range("F6").Value = 0
$wsExport.Range("G6").Value = 8490014.3650100008
$wsExport.Range("H6").Value = 0
$wsExport.Range("I6").Value = 439597.55237699999
$wsExport.Range("J6").Value = "NA"
$wsExport.Range("K6").Value = 0
$wsExport.Range("L6").Value = 0
$wsExport.Range("M6").Value = 20
$wsExport.Range("N6").Value = 10
$wsExport.Range("O6").Value = 0.618176844244679
$wsExport.Range("P6").Value = 0.255757329816374
$wsExport.Range("Q6").Value = 0
$wsExport.Range("R6").Value = 0
$wsExport.Range("S6").Value = "Compound: name of assigned metabolite, noStd means assigment was NOT verified with standard compound"

$wsExport.Activate()
$wsExport.Range("B5").Select() | Out-Null

# Re-activate the "#convert" sheet last, since it is the tab that is
# marked selected/active in the saved workbook.
$wsConvert.Activate()
